$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F12").Value = 0
$ws.Range("F16").Value = -1
